$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item('4_')
$ws5 = $wb.Worksheets.Item('5_')
$ws6 = $wb.Worksheets.Item('6_')

# ---------------------------------------------------------------------------
# Sheet "4_" - Prandtl number lookup practice question
# ---------------------------------------------------------------------------
$ws4.Range('A1').Value = 'What is the Prandtl number for air at 1 atm and 400 degrees C?'
$ws4.Rows.Item(1).RowHeight = 45

$ws4.Range('A2').Value = 0.6948
$ws4.Range('B2').Value = 0.001

$ws4.Range('C3').Value = 'Just a little practice using a material  property table.  '
$ws4.Rows.Item(3).RowHeight = 30

$ws4.Range('C4').Value = 'This is nondimensional, so it should be the same value even if you look at an imperial measurement chart'
$ws4.Rows.Item(4).RowHeight = 60

# ---------------------------------------------------------------------------
# Sheet "5_" - conductivity / Nusselt number conceptual question
# ---------------------------------------------------------------------------
$ws5.Range('A1').Value = 'In the problem from the video, if the thermal conductivity of air was actually higher, what would happen to the Nusselt number?'
$ws5.Rows.Item(1).RowHeight = 90

$ws5.Range('A2').Value = 'It would be smaller'
$ws5.Range('B2').Value = 'Y'
$ws5.Range('C2').Value = 'You can think of this in two ways.  Mathematically, this would make Pr smaller because k is in the numerator of thermal diffusivity, and this would make Nu smaller.  Physically, this would make conduction in the fluid more effective, and so the ratio of convective HT: conductive HT (that is, the Nusselt number) would go down.'
$ws5.Rows.Item(2).RowHeight = 165

$ws5.Range('A3').Value = 'It would stay the same'
$ws5.Range('B3').Value = 'N'

$ws5.Range('A4').Value = 'It would be larger'
$ws5.Range('B4').Value = 'N'

# ---------------------------------------------------------------------------
# Sheet "6_" - Re/Pr -> Nu calculation question
# ---------------------------------------------------------------------------
$ws6.Range('A1').Value = 'Imagine you found the Re number in the example problem to be 40,000, with the same Pr number.  Calculate the Nu number using both of the available coefficients.  How much larger is the Nu using the second equation (for 40K-400K)?'
$ws6.Rows.Item(1).RowHeight = 150

$ws6.Range('A2').Value = 1.76
$ws6.Range('B2').Value = 1
$ws6.Range('C2').Value = 'Two things to note from this: 1) Even though the coefficients look very different, the functions described by them actually meet at the transition points  2) These are mathematical models for experimental data, so they aren''t meant to be perfect (otherwise, the answer to this question would be 0): you must assume a significant error (say about 10%) when using these models.'
$ws6.Rows.Item(2).RowHeight = 180

$ws6.Range('C4').NumberFormat = '#,##0'

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping. "5_" must be activated last so it
# ends up as the workbook's active tab (matches activeTab=5 in the saved
# file), while "4_" and "6_" simply keep their own final selections.
# ---------------------------------------------------------------------------
$ws4.Range('A5').Select()
$ws6.Range('C3').Select()
$ws5.Range('E8').Select()
